$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.419.81"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "1.568.47"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3690"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3310"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.917"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.866"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").Value = "1.567.53"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001110"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06710"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.377"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("E22").Value = "  +1.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").Value = "22.404.14"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.606"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.934"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("D31").Value = "1.748.01"
$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.042"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.976"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.776"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08298"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02425"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.291"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2201"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.261"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6153"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5978"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.17%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.762"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.030"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.89%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07177"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "
